$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1061.3871
$ws.Range("J17").Value = 1061.3871
$ws.Range("L17").Value = 3184.1613
$ws.Range("N17").Value = -3520.1613
$ws.Range("H64").Value = 37205.137
$ws.Range("I64").Value = 74018.42999999999
$ws.Range("J64").Value = 2846.0667
$ws.Range("K64").Value = 74018.42999999999
$ws.Range("L64").Value = 2846.0667
$ws.Range("M64").Value = -73770.42999999999
$ws.Range("N64").Value = -3342.0667
$ws.Range("H67").Value = 37205.137
$ws.Range("I67").Value = 74018.42999999999
$ws.Range("J67").Value = 2846.0667
$ws.Range("K67").Value = 74018.42999999999
$ws.Range("L67").Value = 2846.0667
$ws.Range("M67").Value = -73160.42999999999
$ws.Range("N67").Value = -4562.066699999999
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H80").Value = 26356.59
$ws.Range("I80").Value = 587.1539
$ws.Range("J80").Value = 39241.31
$ws.Range("K80").Value = 1761.4617
$ws.Range("L80").Value = 117723.93
$ws.Range("M80").Value = -763.4617000000001
$ws.Range("N80").Value = -119719.93
$ws.Range("H83").Value = 26356.59
$ws.Range("I83").Value = 587.1539
$ws.Range("J83").Value = 39241.31
$ws.Range("K83").Value = 5284.3851
$ws.Range("L83").Value = 353171.79
$ws.Range("M83").Value = -292.3851000000004
$ws.Range("N83").Value = -363155.79
$ws.Range("H107").Value = 428.21054
$ws.Range("I107").Value = 428.21054
$ws.Range("K107").Value = 428.21054
$ws.Range("M107").Value = 1491.78946
$ws.Range("H112").Value = 1934.091
$ws.Range("J112").Value = 1934.091
$ws.Range("L112").Value = 5802.272999999999
$ws.Range("N112").Value = -8018.272999999999
$ws.Range("H116").Value = 4343.3335
$ws.Range("I116").Value = 2750
$ws.Range("J116").Value = 5140
$ws.Range("K116").Value = 2750
$ws.Range("L116").Value = 5140
$ws.Range("M116").Value = 692
$ws.Range("N116").Value = -12024
$ws.Range("H129").Value = 805.5
$ws.Range("I129").Value = 623.4286
$ws.Range("J129").Value = 819.2043
$ws.Range("K129").Value = 1870.2858
$ws.Range("L129").Value = 2457.6129
$ws.Range("M129").Value = 3129.7142
$ws.Range("N129").Value = -12457.6129

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25945.662
$ws.Range("I32").Value = 5832.1045
$ws.Range("J32").Value = 218461.14
$ws.Range("K32").Value = 5832.1045
$ws.Range("L32").Value = 218461.14
$ws.Range("M32").Value = -5545.1045
$ws.Range("N32").Value = -219035.14
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H61").Value = 2374
$ws.Range("I61").Value = 2311.4
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2311.4
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2099.4
$ws.Range("N61").Value = -3424
$ws.Range("H102").Value = 2845.9092
$ws.Range("I102").Value = 2744.5
$ws.Range("J102").Value = 2903.8572
$ws.Range("K102").Value = 2744.5
$ws.Range("L102").Value = 2903.8572
$ws.Range("M102").Value = -1122.5
$ws.Range("N102").Value = -6147.8572
$ws.Range("H132").Value = 4033
$ws.Range("I132").Value = 4296.25
$ws.Range("J132").Value = 3650.0908
$ws.Range("K132").Value = 12888.75
$ws.Range("L132").Value = 10950.2724
$ws.Range("M132").Value = -10358.75
$ws.Range("N132").Value = -16010.2724
$ws.Range("H136").Value = 2374
$ws.Range("I136").Value = 2311.4
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6934.200000000001
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -4384.200000000001
$ws.Range("N136").Value = -14100

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1780.3334
$ws.Range("I99").Value = 1796.6666
$ws.Range("K99").Value = 1796.6666
$ws.Range("M99").Value = -298.6666
$ws.Range("H105").Value = 2262.8572
$ws.Range("I105").Value = 2472.5
$ws.Range("J105").Value = 1983.3334
$ws.Range("K105").Value = 2472.5
$ws.Range("L105").Value = 1983.3334
$ws.Range("M105").Value = -725.5
$ws.Range("N105").Value = -5477.3334
$ws.Range("H134").Value = 2196.1035
$ws.Range("I134").Value = 2159.75
$ws.Range("K134").Value = 6479.25
$ws.Range("M134").Value = -3944.25

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 532.8182
$ws.Range("I16").Value = 441.33334
$ws.Range("J16").Value = 642.6
$ws.Range("K16").Value = 441.33334
$ws.Range("L16").Value = 642.6
$ws.Range("M16").Value = -154.33334
$ws.Range("N16").Value = -1216.6
$ws.Range("H99").Value = 7186.091
$ws.Range("I99").Value = 3231.5557
$ws.Range("J99").Value = 9923.846
$ws.Range("K99").Value = 3231.5557
$ws.Range("L99").Value = 9923.846
$ws.Range("M99").Value = -1733.5557
$ws.Range("N99").Value = -12919.846
$ws.Range("H105").Value = 1010.6818
$ws.Range("I105").Value = 875.8889
$ws.Range("K105").Value = 875.8889
$ws.Range("M105").Value = 871.1111
$ws.Range("H107").Value = 678.2
$ws.Range("I107").Value = 745.6
$ws.Range("K107").Value = 745.6
$ws.Range("M107").Value = 1174.4
$ws.Range("H113").Value = 532.8182
$ws.Range("I113").Value = 441.33334
$ws.Range("J113").Value = 642.6
$ws.Range("K113").Value = 441.33334
$ws.Range("L113").Value = 642.6
$ws.Range("M113").Value = 1728.66666
$ws.Range("N113").Value = -4982.6
$ws.Range("H122").Value = 530.2727
$ws.Range("I122").Value = 296.66666
$ws.Range("J122").Value = 617.875
$ws.Range("K122").Value = 889.9999799999999
$ws.Range("L122").Value = 1853.625
$ws.Range("M122").Value = 1560.00002
$ws.Range("N122").Value = -6753.625
$ws.Range("H126").Value = 7186.091
$ws.Range("I126").Value = 3231.5557
$ws.Range("J126").Value = 9923.846
$ws.Range("K126").Value = 9694.667099999999
$ws.Range("L126").Value = 29771.538
$ws.Range("M126").Value = -7224.667099999999
$ws.Range("N126").Value = -34711.538

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1388
$ws.Range("N4").ClearContents()
$ws.Range("H14").Value = 270.3846
$ws.Range("I14").Value = 270.3846
$ws.Range("K14").Value = 811.1537999999999
$ws.Range("M14").Value = -638.1537999999999
$ws.Range("H115").Value = 3435.0557
$ws.Range("I115").Value = 824.6667
$ws.Range("J115").Value = 3957.1333
$ws.Range("K115").Value = 2474.0001
$ws.Range("L115").Value = 11871.3999
$ws.Range("M115").Value = -1299.0001
$ws.Range("N115").Value = -14221.3999
$ws.Range("H131").Value = 6787.3765
$ws.Range("I131").Value = 1567.5
$ws.Range("J131").Value = 7021.9775
$ws.Range("K131").Value = 4702.5
$ws.Range("L131").Value = 21065.9325
$ws.Range("M131").Value = 337.5
$ws.Range("N131").Value = -31145.9325

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 125001100
$ws.Range("I97").Value = 125001100
$ws.Range("K97").Value = 125001100
$ws.Range("M97").Value = -125000604
$ws.Range("H102").Value = 2748.4375
$ws.Range("I102").Value = 2338.3333
$ws.Range("J102").Value = 2994.5
$ws.Range("K102").Value = 2338.3333
$ws.Range("L102").Value = 2994.5
$ws.Range("M102").Value = -716.3332999999998
$ws.Range("N102").Value = -6238.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1309.4138
$ws.Range("I22").Value = 1870.3
$ws.Range("J22").Value = 1014.2105
$ws.Range("K22").Value = 1870.3
$ws.Range("L22").Value = 1014.2105
$ws.Range("M22").Value = -1575.3
$ws.Range("N22").Value = -1604.2105
$ws.Range("H27").Value = 1309.4138
$ws.Range("I27").Value = 1870.3
$ws.Range("J27").Value = 1014.2105
$ws.Range("K27").Value = 1870.3
$ws.Range("L27").Value = 1014.2105
$ws.Range("M27").Value = -1763.3
$ws.Range("N27").Value = -1228.2105
$ws.Range("H61").Value = 1583.9231
$ws.Range("I61").Value = 1569.2307
$ws.Range("J61").Value = 1598.6154
$ws.Range("K61").Value = 1569.2307
$ws.Range("L61").Value = 1598.6154
$ws.Range("M61").Value = -1367.2307
$ws.Range("N61").Value = -2002.6154
$ws.Range("H93").Value = 1441.8
$ws.Range("I93").Value = 1668.7059
$ws.Range("K93").Value = 1668.7059
$ws.Range("M93").Value = -420.7058999999999
$ws.Range("H100").Value = 2196.5557
$ws.Range("I100").Value = 1781.5
$ws.Range("J100").Value = 3026.6667
$ws.Range("K100").Value = 1781.5
$ws.Range("L100").Value = 3026.6667
$ws.Range("M100").Value = -1240.5
$ws.Range("N100").Value = -4108.6667
$ws.Range("H113").Value = 1583.9231
$ws.Range("I113").Value = 1569.2307
$ws.Range("J113").Value = 1598.6154
$ws.Range("K113").Value = 1569.2307
$ws.Range("L113").Value = 1598.6154
$ws.Range("M113").Value = 600.7692999999999
$ws.Range("N113").Value = -5938.6154
$ws.Range("H119").Value = 30283.334
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 30283.334
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 30283.334
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -39959.334
$ws.Range("H136").Value = 1374.5366
$ws.Range("I136").Value = 1225.1765
$ws.Range("J136").Value = 2100
$ws.Range("K136").Value = 3675.5295
$ws.Range("L136").Value = 6300
$ws.Range("M136").Value = -1125.5295
$ws.Range("N136").Value = -11400

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 41881.332
$ws.Range("J119").Value = 41881.332
$ws.Range("L119").Value = 41881.332
$ws.Range("N119").Value = -51557.332
$ws.Range("H122").Value = 1598
$ws.Range("I122").Value = 1397.1428
$ws.Range("J122").Value = 1773.75
$ws.Range("K122").Value = 4191.428400000001
$ws.Range("L122").Value = 5321.25
$ws.Range("M122").Value = -1741.428400000001
$ws.Range("N122").Value = -10221.25
$ws.Range("H126").Value = 1008.5789
$ws.Range("I126").Value = 885.8125
$ws.Range("J126").Value = 1663.3334
$ws.Range("K126").Value = 2657.4375
$ws.Range("L126").Value = 4990.0002
$ws.Range("M126").Value = -187.4375
$ws.Range("N126").Value = -9930.0002
$ws.Range("H132").Value = 2031.2609
$ws.Range("I132").Value = 2074.4866
$ws.Range("J132").Value = 1853.5555
$ws.Range("K132").Value = 6223.459800000001
$ws.Range("L132").Value = 5560.666499999999
$ws.Range("M132").Value = -3693.459800000001
$ws.Range("N132").Value = -10620.6665
$ws.Range("H135").Value = 39900
$ws.Range("J135").Value = 39900
$ws.Range("L135").Value = 39900
$ws.Range("N135").Value = -50040
$ws.Range("H136").Value = 3838.8
$ws.Range("I136").Value = 923
$ws.Range("K136").Value = 2769
$ws.Range("M136").Value = -219

